$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# "añadido moto2 y moto3 a 2013": the Hoja1 lookup row that pointed at the
# 2012 season workbook now points at the 2013 one (which has had its Moto2 /
# Moto3 sheets added), so the file_path / file_name pair moves from 2012 to
# 2013.
$ws.Range("A2").Value = "C:\Users\zaka\Desktop\MOTOGP\Excels\data\2013.xlsx"
$ws.Range("B2").Value = "2013"

# The saved workbook also recorded the cursor resting on B5 afterwards.
$ws.Range("B5").Select()
